# Updated cryptos list on Sun Aug 18 16:56:21 UTC 2024 with GitHub Actions
#
# This script updates the "Price" (D) and "Volume(1h)" (E) columns for the
# crypto ranking sheet, and fixes a couple of rows whose Coin/Link/Price/
# Volume values were swapped (Stacks <-> Fetch.AI, Maker <-> RenderToken).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value as plain text, even if it looks like a number
# (e.g. "21.26"), without leaving a lingering number-format override on
# the cell once we're done.
function Set-TextValue($addr, $value) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $value
    $r.ClearFormats()
}

# Row 2 - Bitcoin
Set-TextValue "D2" "59.834.11"
$ws.Range("E2").Value = "  +0.66%  "

# Row 3 - Ethereum
Set-TextValue "D3" "2.654.11"
$ws.Range("E3").Value = "  +1.84%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.04%  "

# Row 5 - BNB
Set-TextValue "D5" "538.75"
$ws.Range("E5").Value = "  +0.57%  "

# Row 6 - Solana
Set-TextValue "D6" "146.32"
$ws.Range("E6").Value = "  +3.81%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  +0.00%  "

# Row 8 - XRP
Set-TextValue "D8" "0.575"
$ws.Range("E8").Value = "  +1.25%  "

# Row 9 - LidoStakedEther
Set-TextValue "D9" "2.673.57"
$ws.Range("E9").Value = "  +2.39%  "

# Row 10 - Toncoin
$ws.Range("E10").Value = "  +3.05%  "

# Row 11 - Dogecoin
Set-TextValue "D11" "0.104"
$ws.Range("E11").Value = "  +0.62%  "

# Row 12 - Cardano
$ws.Range("E12").Value = "  +0.76%  "

# Row 13 - TRON
$ws.Range("E13").Value = "  -0.59%  "

# Row 14 - WrappedliquidstakedEther2.0
Set-TextValue "D14" "3.125.97"
$ws.Range("E14").Value = "  +2.04%  "

# Row 15 - WrappedBTC
Set-TextValue "D15" "59.752.16"
$ws.Range("E15").Value = "  +0.66%  "

# Row 16 - Avalanche
Set-TextValue "D16" "21.26"
$ws.Range("E16").Value = "  +3.33%  "

# Row 17 - WrappedEther
Set-TextValue "D17" "2.673.51"
$ws.Range("E17").Value = "  +1.62%  "

# Row 18 - ShibaInu
$ws.Range("E18").Value = "  +1.44%  "

# Row 19 - BitcoinCash
$ws.Range("E19").Value = "  -0.73%  "

# Row 20 - Polkadot
$ws.Range("E20").Value = "  +1.96%  "

# Row 21 - Chainlink
Set-TextValue "D21" "10.42"
$ws.Range("E21").Value = "  +2.93%  "

# Row 22 - Uniswap
$ws.Range("E22").Value = "  -0.28%  "

# Row 23 - Dai
$ws.Range("E23").Value = "  +0.03%  "

# Row 24 - Litecoin
Set-TextValue "D24" "66.71"
$ws.Range("E24").Value = "  -0.50%  "

# Row 25 - Polygon
$ws.Range("E25").Value = "  +2.13%  "

# Row 26 - Kaspa
$ws.Range("E26").Value = "  -0.59%  "

# Row 27 - Binance-PegBSC-USD
$ws.Range("E27").Value = "  -0.09%  "

# Row 28 - InternetComputer(DFINITY)
$ws.Range("E28").Value = "  +2.18%  "

# Row 29 - PEPE
$ws.Range("E29").Value = "  +1.89%  "

# Row 30 - USDe
$ws.Range("E30").Value = "  -0.08%  "

# Row 31 - PancakeSwap
$ws.Range("E31").Value = "  +1.35%  "

# Row 32 - Aptos
Set-TextValue "D32" "5.87"
$ws.Range("E32").Value = "  +0.51%  "

# Row 33 - EthereumClassic
Set-TextValue "D33" "19.01"
$ws.Range("E33").Value = "  +0.90%  "

# Row 34 - Monero
Set-TextValue "D34" "150.68"
$ws.Range("E34").Value = "  +0.91%  "

# Row 35 - NEARProtocol
$ws.Range("E35").Value = "  +1.16%  "

# Row 36 - ImmutableX
$ws.Range("E36").Value = "  +1.93%  "

# Row 37 - SuiNetwork
Set-TextValue "D37" "0.840"
$ws.Range("E37").Value = "  -0.76%  "

# Row 38 / Row 39 - Stacks and Fetch.AI swap places
Set-TextValue "B38" "Fetch.AI"
Set-TextValue "C38" "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
Set-TextValue "D38" "0.845"
$ws.Range("E38").Value = "  +0.78%  "

Set-TextValue "B39" "Stacks"
Set-TextValue "C39" "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
Set-TextValue "D39" "1.47"
$ws.Range("E39").Value = "  -0.18%  "

# Row 40 - Bittensor
Set-TextValue "D40" "292.22"
$ws.Range("E40").Value = "  +5.14%  "

# Row 41 - Filecoin
$ws.Range("E41").Value = "  +1.84%  "

# Row 42 - FirstDigitalUSD
Set-TextValue "D42" "0.999"
$ws.Range("E42").Value = "  +0.11%  "

# Row 43 - Mantle
Set-TextValue "D43" "0.610"
$ws.Range("E43").Value = "  +1.71%  "

# Row 44 - EnergySwap
Set-TextValue "D44" "19.52"
$ws.Range("E44").Value = "  +4.91%  "

# Row 45 - Hedera
Set-TextValue "D45" "0.0540"

# Row 46 - WhiteBITCoin
Set-TextValue "D46" "10.72"
$ws.Range("E46").Value = "  -0.42%  "

# Row 47 - Stellar
$ws.Range("E47").Value = "  -1.04%  "

# Row 48 / Row 50 - Maker and RenderToken swap places
Set-TextValue "B48" "RenderToken"
Set-TextValue "C48" "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue "D48" "4.75"
$ws.Range("E48").Value = "  +4.90%  "

# Row 49 - VeChain
$ws.Range("E49").Value = "  +2.01%  "

Set-TextValue "B50" "Maker"
Set-TextValue "C50" "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
Set-TextValue "D50" "1.980.69"
$ws.Range("E50").Value = "  +1.51%  "

# Row 51 - InjectiveProtocol
Set-TextValue "D51" "18.46"
$ws.Range("E51").Value = "  +0.46%  "
